$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: I1 = "I0", J1 = "IF" (same style as other headers, e.g. H1)
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-48
$iVals = @(6,9,9,8,8,9,10,7,8,8,10,6,8,8,8,9,9,9,7,7,8,8,5,8,9,6,9,4,7,9,9,8,9,9,6,8,6,8,9,9,6,6,6,5,6,3,5)
$jVals = @(7,9,9,9,8,9,10,7,8,8,10,6,8,8,8,9,9,9,8,7,8,8,7,8,9,7,9,4,9,9,9,8,9,9,7,8,6,8,9,9,6,6,6,5,6,3,5)

for ($idx = 0; $idx -lt $iVals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$idx]
    $ws.Cells.Item($row, 10).Value = $jVals[$idx]
}
